# Update the build timestamp embedded in the workbook's version strings
# from "February 03 2026 17.29.55 EST" to "February 03 2026 18.05.36 EST".

$wb = $excel.ActiveWorkbook

$newStamp = "February 03 2026 18.05.36 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet ---

$aboutSheet.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

# Build the citation text with a single-quoted literal (so the embedded
# double-quote and single-quote characters need no backtick escaping), then
# splice in the computed timestamp.
$citationPrefix = 'Recommended Citation:  "Global Energy Monitor, Coal mine boundaries and methane sources for Zhujidong Coal Mine, China, M1254, version ''Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on '
$citationSuffix = ')''. (See the CC license for attribution requirements if sharing or adapting the data set.)'

$aboutSheet.Range("A6").Value = $citationPrefix + $newStamp + $citationSuffix

# --- "Boundaries and methane sources" sheet ---

$newVersionText = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

foreach ($r in 2..7) {
    $dataSheet.Range("S$r").Value = $newVersionText
}
